$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text columns (Coin name / Link URL): plain text, no numeric coercion risk ---
$ws.Range("B6").Value = "GateToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("B8").Value = "KuCoinToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("B9").Value = "BTSEToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("B10").Value = "MXToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("B12").Value = "WazirX"
$ws.Range("C12").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("B13").Value = "MandalaExchangeToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"

# --- Numeric-looking columns (Price / Volume%): force Text format so Excel keeps
# the original string (with trailing zeros / percent sign) instead of coercing to a number ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "313.25"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "37.69"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "1.44%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.132"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.55%"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "2.10%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.412"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.49%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.930"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "1.64%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "8.285"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "0.86%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.919"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-4.61%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9220"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "0.15%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1230"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-1.84%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1931"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "2.05%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09143"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "4.37%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03310"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-3.36%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09632"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.86%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001378"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.50%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.005734"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-5.12%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.511"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-1.49%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "2.15%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.250"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "4.21%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1272"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-1.11%"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "3.33%"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.71%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.04369"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.46%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001249"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "1.98%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004312"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-3.94%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001220"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-10.15%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02220"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "2.32%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05123"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "4.04%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007464"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-2.92%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1365"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "2.41%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.008781"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-10.85%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.001959"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-2.27%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.008612"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-2.45%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006727"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-1.67%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.57%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003344"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "10.87%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.001200"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-8.23%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002100"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.57%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002000"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.57%"
